$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (CD45 / Identity / X) is a duplicate of row 2 (CD45) and is removed,
# shifting all subsequent rows up by one.
$ws.Rows(14).Delete()

# Update the Cell_Type column (C) values for the Identity markers that keep
# an "X" marking to the lower-case "x", and fill in the previously blank
# Cell_Type cell for CD45 (row 2) with "x" as well.
$ws.Range("C2").Value = "x"
$ws.Range("C13").Value = "x"
$ws.Range("C14").Value = "x"
